$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to stay text so values like "1.003" or
# "  -0.69%  " are not silently reinterpreted as numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "26.037.31"
$ws.Cells.Item(2, 5).Value = "  -0.69%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.649.10"
$ws.Cells.Item(3, 5).Value = "  -0.62%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "1.003"
$ws.Cells.Item(4, 5).Value = "  -0.17%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "217.93"
$ws.Cells.Item(5, 5).Value = "  -0.42%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "0.5209"
$ws.Cells.Item(6, 5).Value = "  -0.26%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.08%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.2612"
$ws.Cells.Item(8, 5).Value = "  -2.02%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.06263"
$ws.Cells.Item(9, 5).Value = "  -1.09%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "20.42"
$ws.Cells.Item(10, 5).Value = "  -3.74%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "0.07737"
$ws.Cells.Item(11, 5).Value = "  -0.09%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "Polkadot"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(12, 4).Value = "4.458"
$ws.Cells.Item(12, 5).Value = "  +0.41%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.622.78"
$ws.Cells.Item(13, 5).Value = "  -2.11%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "0.5426"
$ws.Cells.Item(14, 5).Value = "  -1.10%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "0.0₅8073"
$ws.Cells.Item(15, 5).Value = "  -2.29%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "64.68"
$ws.Cells.Item(16, 5).Value = "  -0.52%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "26.062.68"
$ws.Cells.Item(17, 5).Value = "  -0.61%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "1.002"
$ws.Cells.Item(18, 5).Value = "  -0.29%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "4.562"
$ws.Cells.Item(19, 5).Value = "  -2.83%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "191.77"
$ws.Cells.Item(20, 5).Value = "  -0.74%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "10.01"
$ws.Cells.Item(21, 5).Value = "  -1.95%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "5.977"
$ws.Cells.Item(22, 5).Value = "  -2.55%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "1.004"
$ws.Cells.Item(23, 5).Value = "  -0.22%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "138.60"
$ws.Cells.Item(24, 5).Value = "  -0.04%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "0.1230"
$ws.Cells.Item(25, 5).Value = "  -0.77%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "7.240"
$ws.Cells.Item(26, 5).Value = "  -0.60%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +0.02%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "1.398"
$ws.Cells.Item(28, 5).Value = "  -1.24%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "0.05919"
$ws.Cells.Item(29, 5).Value = "  -2.32%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "1.274"
$ws.Cells.Item(30, 5).Value = "  -0.81%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "3.491"
$ws.Cells.Item(31, 5).Value = "  -1.97%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "3.224"
$ws.Cells.Item(32, 5).Value = "  -3.94%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "1.523"
$ws.Cells.Item(33, 5).Value = "  -8.06%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "HuobiToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(34, 4).Value = "2.414"
$ws.Cells.Item(34, 5).Value = "  +0.13%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "ARBITRUM"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(35, 4).Value = "0.9445"
$ws.Cells.Item(35, 5).Value = "  -4.22%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.97%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "0.5738"
$ws.Cells.Item(37, 5).Value = "  -3.44%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.01599"
$ws.Cells.Item(38, 5).Value = "  +0.05%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "5.850"
$ws.Cells.Item(39, 5).Value = "  -2.04%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "0.8446"
$ws.Cells.Item(40, 5).Value = "  -2.53%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "1.002"
$ws.Cells.Item(41, 5).Value = "  -0.14%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "100.50"
$ws.Cells.Item(42, 5).Value = "  +0.51%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "1.000.36"
$ws.Cells.Item(43, 5).Value = "  -4.65%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "1.794.03"
$ws.Cells.Item(44, 5).Value = "  -0.15%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "56.53"
$ws.Cells.Item(45, 5).Value = "  -1.54%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "0.0₈105"
$ws.Cells.Item(46, 5).Value = "  -4.12%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "1.002"
$ws.Cells.Item(47, 5).Value = "  -0.48%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "0.4292"
$ws.Cells.Item(48, 5).Value = "  +1.51%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "7.846"
$ws.Cells.Item(49, 5).Value = "  -3.58%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -0.63%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "1.463"
$ws.Cells.Item(51, 5).Value = "  -0.99%  "
